$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.227.63"
$ws.Range("E2").Value2 = "  -0.79%  "
$ws.Range("D3").Value2 = "3.492.19"
$ws.Range("E3").Value2 = "  -0.19%  "
$ws.Range("E4").Value2 = "  -0.06%  "
$ws.Range("D5").Value2 = "'604.36"
$ws.Range("E5").Value2 = "  +0.54%  "
$ws.Range("D6").Value2 = "'144.02"
$ws.Range("E6").Value2 = "  -2.37%  "
$ws.Range("D7").Value2 = "3.491.77"
$ws.Range("E7").Value2 = "  -0.21%  "
$ws.Range("E8").Value2 = "  -0.07%  "
$ws.Range("E9").Value2 = "  -0.66%  "
$ws.Range("D10").Value2 = "'8.05"
$ws.Range("E10").Value2 = "  +1.61%  "
$ws.Range("E11").Value2 = "  -4.65%  "
$ws.Range("D12").Value2 = "'0.412"
$ws.Range("E12").Value2 = "  -2.46%  "
$ws.Range("D13").Value2 = "4.081.49"
$ws.Range("E13").Value2 = "  -0.26%  "
$ws.Range("B14").Value2 = "ShibaInu"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value2 = "'0.0000203"
$ws.Range("E14").Value2 = "  -4.69%  "
$ws.Range("B15").Value2 = "Avalanche"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value2 = "'30.34"
$ws.Range("E15").Value2 = "  -2.69%  "
$ws.Range("D16").Value2 = "3.489.51"
$ws.Range("E16").Value2 = "  -0.32%  "
$ws.Range("D17").Value2 = "66.256.68"
$ws.Range("E17").Value2 = "  -0.79%  "
$ws.Range("D18").Value2 = "'0.117"
$ws.Range("E18").Value2 = "  -0.35%  "
$ws.Range("E19").Value2 = "  +3.42%  "
$ws.Range("D20").Value2 = "'6.15"
$ws.Range("E20").Value2 = "  -3.69%  "
$ws.Range("D21").Value2 = "'14.84"
$ws.Range("E21").Value2 = "  -3.33%  "
$ws.Range("D22").Value2 = "'425.69"
$ws.Range("E22").Value2 = "  -1.81%  "
$ws.Range("D23").Value2 = "'0.593"
$ws.Range("E23").Value2 = "  -2.67%  "
$ws.Range("D24").Value2 = "'77.92"
$ws.Range("E24").Value2 = "  -2.10%  "
$ws.Range("D25").Value2 = "3.622.03"
$ws.Range("E25").Value2 = "  -0.46%  "
$ws.Range("E26").Value2 = "  +0.11%  "
$ws.Range("D27").Value2 = "'0.0000117"
$ws.Range("E27").Value2 = "  -1.95%  "
$ws.Range("D28").Value2 = "'9.28"
$ws.Range("E28").Value2 = "  -5.71%  "
$ws.Range("E29").Value2 = "  -4.05%  "
$ws.Range("D30").Value2 = "'2.46"
$ws.Range("E30").Value2 = "  -1.04%  "
$ws.Range("D31").Value2 = "'1.01"
$ws.Range("E31").Value2 = "  +1.44%  "
$ws.Range("E32").Value2 = "  -0.14%  "
$ws.Range("E33").Value2 = "  -8.77%  "
$ws.Range("D34").Value2 = "'25.11"
$ws.Range("E34").Value2 = "  -1.11%  "
$ws.Range("D35").Value2 = "3.476.95"
$ws.Range("E35").Value2 = "  -0.48%  "
$ws.Range("D37").Value2 = "'1.74"
$ws.Range("E37").Value2 = "  -3.41%  "
$ws.Range("D38").Value2 = "'5.62"
$ws.Range("E38").Value2 = "  -4.96%  "
$ws.Range("D39").Value2 = "'7.72"
$ws.Range("E39").Value2 = "  -3.42%  "
$ws.Range("E40").Value2 = "  -0.01%  "
$ws.Range("D41").Value2 = "'169.94"
$ws.Range("E41").Value2 = "  -0.14%  "
$ws.Range("D42").Value2 = "'0.0858"
$ws.Range("E42").Value2 = "  -3.82%  "
$ws.Range("E43").Value2 = "  -4.96%  "
$ws.Range("D44").Value2 = "'0.881"
$ws.Range("E44").Value2 = "  -1.72%  "
$ws.Range("E45").Value2 = "  -8.71%  "
$ws.Range("D46").Value2 = "'45.45"
$ws.Range("E46").Value2 = "  -0.85%  "
$ws.Range("D47").Value2 = "'25.88"
$ws.Range("E47").Value2 = "  -9.52%  "
$ws.Range("E48").Value2 = "  -10.40%  "
$ws.Range("D49").Value2 = "'2.41"
$ws.Range("E49").Value2 = "  -0.72%  "
$ws.Range("D50").Value2 = "'7.13"
$ws.Range("E50").Value2 = "  -4.49%  "
$ws.Range("E51").Value2 = "  -2.94%  "
